{"js": "// Update the date line and each division problem in the table to the new values.\nconst replacements = [\n  [\"2024-09-16 Monday\", \"2024-09-17 Tuesday\"],\n  [\"780\u00f72=\", \"194\u00f77=\"],\n  [\"477\u00f72=\", \"551\u00f78=\"],\n  [\"415\u00f78=\", \"708\u00f79=\"],\n  [\"271\u00f72=\", \"247\u00f76=\"],\n  [\"529\u00f72=\", \"233\u00f77=\"],\n  [\"201\u00f79=\", \"304\u00f79=\"],\n  [\"310\u00f77=\", \"259\u00f74=\"],\n  [\"911\u00f78=\", \"408\u00f73=\"],\n  [\"413\u00f78=\", \"313\u00f77=\"],\n  [\"947\u00f73=\", \"148\u00f74=\"],\n  [\"848\u00f79=\", \"654\u00f75=\"],\n  [\"682\u00f74=\", \"112\u00f76=\"],\n  [\"863\u00f77=\", \"550\u00f72=\"],\n  [\"970\u00f75=\", \"853\u00f75=\"],\n  [\"668\u00f76=\", \"613\u00f73=\"],\n  [\"168\u00f73=\", \"845\u00f75=\"],\n  [\"832\u00f76=\", \"653\u00f79=\"],\n  [\"727\u00f76=\", \"276\u00f76=\"],\n  [\"149\u00f76=\", \"253\u00f79=\"],\n  [\"337\u00f74=\", \"732\u00f77=\"],\n  [\"384\u00f75=\", \"263\u00f79=\"],\n  [\"912\u00f72=\", \"979\u00f74=\"],\n  [\"530\u00f77=\", \"397\u00f79=\"],\n  [\"907\u00f73=\", \"785\u00f74=\"],\n  [\"855\u00f77=\", \"224\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each division problem in the table to the new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-09-16 Monday', '2024-09-17 Tuesday'),\n    @('780\u00f72=', '194\u00f77='),\n    @('477\u00f72=', '551\u00f78='),\n    @('415\u00f78=', '708\u00f79='),\n    @('271\u00f72=', '247\u00f76='),\n    @('529\u00f72=', '233\u00f77='),\n    @('201\u00f79=', '304\u00f79='),\n    @('310\u00f77=', '259\u00f74='),\n    @('911\u00f78=', '408\u00f73='),\n    @('413\u00f78=', '313\u00f77='),\n    @('947\u00f73=', '148\u00f74='),\n    @('848\u00f79=', '654\u00f75='),\n    @('682\u00f74=', '112\u00f76='),\n    @('863\u00f77=', '550\u00f72='),\n    @('970\u00f75=', '853\u00f75='),\n    @('668\u00f76=', '613\u00f73='),\n    @('168\u00f73=', '845\u00f75='),\n    @('832\u00f76=', '653\u00f79='),\n    @('727\u00f76=', '276\u00f76='),\n    @('149\u00f76=', '253\u00f79='),\n    @('337\u00f74=', '732\u00f77='),\n    @('384\u00f75=', '263\u00f79='),\n    @('912\u00f72=', '979\u00f74='),\n    @('530\u00f77=', '397\u00f79='),\n    @('907\u00f73=', '785\u00f74='),\n    @('855\u00f77=', '224\u00f76='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Search text not found: $oldText\"\n    }\n}\n"}
